$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previously used range so stale cells (B1, C1, D1, F1, G1) are removed
$ws.Cells.Clear()

# Set the new values per the diff
$ws.Range("A1").Value = "Q5"
$ws.Range("B2").Value = "VORNE"
